# "merge Huy, Thong files" - update Sheet1!E5 from "HashMap" to "Trie tree"
# (C5 keeps its existing text "Hashtable, ADT, Analysis of Algorithms")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E5").Value = "Trie tree"

# Match the author's recorded view/selection state on save: scrolled so
# row 4 is the top visible row, column B is the left visible column, and
# E5 is the selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("E5").Select()
